$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update "想去人数" (want-to-go count) column F
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1081
$wsExhibit.Range("F4").Value = 1647
$wsExhibit.Range("F5").Value = 750
$wsExhibit.Range("F6").Value = 153

# Sheet "全部类型" (all types) - update "想去人数" (want-to-go count) column F
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1081
$wsAll.Range("F4").Value = 1647
$wsAll.Range("F6").Value = 750
$wsAll.Range("F7").Value = 153
